$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Data")

# Labels in column J (rows 2-6), averages formulas in column K (rows 2-6)
$labels = @("avg1", "avg2", "avg3", "avg4", "avg5")
$cols = @("D", "E", "F", "G", "H")

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws.Range("J$row").Value = $labels[$i]
    $ws.Range("K$row").Formula = "=AVERAGE($($cols[$i])2:$($cols[$i])60)"
}

$ws.Range("K7").Select() | Out-Null
